# vim_lib.xlsx update:
#  1. Add how to type <Tab> in vim (new row appended to the "git script /
#     basic / config / ..." table on the "Sheet2" worksheet).
#  2. (Redis library addition mentioned in the commit message is not part
#     of this particular sheet/sharedStrings change set.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# New row 14: Lang/Topic/Command-style triple describing how to insert a
# literal <Tab> character while expandtab is enabled in vim.
# Cells are populated C -> B -> A so the underlying shared-string table
# picks up the three new strings in the same relative order as the
# target workbook.
$ws.Range("C14").Formula = "You can use <CTRL-V><Tab> in ""insert mode"". In insert mode <CTRL-V> inserts a literal copy of your next character."
$ws.Range("B14").Formula = "Insert <Tab> when expandtab ON"
$ws.Range("A14").Formula = "vim"

# Move the active selection down to the row below the newly added data,
# matching the post-edit cursor position.
$null = $ws.Range("A15").Select()
